# Update Assignment 2 grades (column C) for re-graded students.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value  = 60.64776531782864
$ws.Range("C5").Value  = 67.02237126211561
$ws.Range("C6").Value  = 67.02237126211561
$ws.Range("C7").Value  = 84.66243092374654
$ws.Range("C8").Value  = 32.96269956291068
$ws.Range("C10").Value = 55.3506518718915
$ws.Range("C11").Value = 15.04705546591035
$ws.Range("C12").Value = 31.85099496236071
$ws.Range("C14").Value = 55.3506518718915
$ws.Range("C15").Value = 21.85298835372897
$ws.Range("C16").Value = 23.62611504402493
$ws.Range("C20").Value = 23.62611504402493
$ws.Range("C21").Value = 31.85099496236071
$ws.Range("C22").Value = 31.85099496236071
$ws.Range("C23").Value = 21.85298835372897
$ws.Range("C24").Value = 84.66243092374654
